$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 0.1952936666666667
$ws.Range("H2").Value = 0.585881
$ws.Range("I2").Value = 0.001827617096392301
$ws.Range("J2").Value = 0.0018276170963923
$ws.Range("M2").Value = 15.090721
$ws.Range("N2").Value = 45.272163
$ws.Range("O2").Value = 0.169971412714946
$ws.Range("P2").Value = 0.169971412714946
$ws.Range("Q2").Value = 2.947122236733667
$ws.Range("R2").Value = 26.524100130603
$ws.Range("S2").Value = 0.000310642659775787
$ws.Range("T2").Value = 0.0003106426597757869

$ws.Range("G3").Value = 0.1952936666666667
$ws.Range("H3").Value = 0.585881
$ws.Range("I3").Value = 0.001827617096392301
$ws.Range("J3").Value = 0.0018276170963923
$ws.Range("N3").Value = 60.45961299999999
$ws.Range("O3").Value = 0.226991713071207
$ws.Range("P3").Value = 0.226991713071207
$ws.Range("Q3").Value = 3.935793169339222
$ws.Range("R3").Value = 35.422138524053
$ws.Range("S3").Value = 0.0004148539355483135
$ws.Range("T3").Value = 0.0004148539355483134

$ws.Range("G4").Value = 0.1952936666666667
$ws.Range("H4").Value = 0.585881
$ws.Range("I4").Value = 0.001827617096392301
$ws.Range("J4").Value = 0.0018276170963923
$ws.Range("M4").Value = 18.66868666666667
$ws.Range("N4").Value = 56.00606
$ws.Range("O4").Value = 0.2102711359030499
$ws.Range("P4").Value = 0.2102711359030499
$ws.Range("Q4").Value = 3.645876270984445
$ws.Range("R4").Value = 32.81288643886
$ws.Range("S4").Value = 0.0003842951228542429
$ws.Range("T4").Value = 0.0003842951228542428

$ws.Range("G5").Value = 0.1952936666666667
$ws.Range("H5").Value = 0.585881
$ws.Range("I5").Value = 0.001827617096392301
$ws.Range("J5").Value = 0.0018276170963923
$ws.Range("M5").Value = 5.641943333333334
$ws.Range("N5").Value = 16.92583
$ws.Range("O5").Value = 0.06354693581733691
$ws.Range("P5").Value = 0.06354693581733692
$ws.Range("Q5").Value = 1.101835800692222
$ws.Range("R5").Value = 9.916522206230001
$ws.Range("S5").Value = 0.0001161394663231092
$ws.Range("T5").Value = 0.0001161394663231091

$ws.Range("G6").Value = 0.1952936666666667
$ws.Range("H6").Value = 0.585881
$ws.Range("I6").Value = 0.001827617096392301
$ws.Range("J6").Value = 0.0018276170963923
$ws.Range("M6").Value = 15.42507366666667
$ws.Range("N6").Value = 46.275221
$ws.Range("O6").Value = 0.1737373292074942
$ws.Range("P6").Value = 0.1737373292074942
$ws.Range("Q6").Value = 3.012419194966778
$ws.Range("R6").Value = 27.111772754701
$ws.Range("S6").Value = 0.0003175253131411537
$ws.Range("T6").Value = 0.0003175253131411537

$ws.Range("G7").Value = 0.1952936666666667
$ws.Range("H7").Value = 0.585881
$ws.Range("I7").Value = 0.001827617096392301
$ws.Range("J7").Value = 0.0018276170963923
$ws.Range("M7").Value = 13.804248
$ws.Range("N7").Value = 41.412744
$ws.Range("O7").Value = 0.1554814732859661
$ws.Range("P7").Value = 0.1554814732859661
$ws.Range("Q7").Value = 2.695882207496
$ws.Range("R7").Value = 24.262939867464
$ws.Range("S7").Value = 0.0002841605987496945
$ws.Range("T7").Value = 0.0002841605987496944

$ws.Range("G8").Value = 3.363724333333333
$ws.Range("H8").Value = 10.091173
$ws.Range("I8").Value = 0.03147874789838274
$ws.Range("J8").Value = 0.03147874789838274
$ws.Range("M8").Value = 15.090721
$ws.Range("N8").Value = 45.272163
$ws.Range("O8").Value = 0.169971412714946
$ws.Range("P8").Value = 0.169971412714946
$ws.Range("Q8").Value = 50.76102543524433
$ws.Range("R8").Value = 456.849228917199
$ws.Range("S8").Value = 0.005350487250785753
$ws.Range("T8").Value = 0.005350487250785752

$ws.Range("G9").Value = 3.363724333333333
$ws.Range("H9").Value = 10.091173
$ws.Range("I9").Value = 0.03147874789838274
$ws.Range("J9").Value = 0.03147874789838274
$ws.Range("N9").Value = 60.45961299999999
$ws.Range("O9").Value = 0.226991713071207
$ws.Range("P9").Value = 0.226991713071207
$ws.Range("Q9").Value = 67.7898238106721
$ws.Range("R9").Value = 610.1084142960489
$ws.Range("S9").Value = 0.007145414910790556
$ws.Range("T9").Value = 0.007145414910790555

$ws.Range("G10").Value = 3.363724333333333
$ws.Range("H10").Value = 10.091173
$ws.Range("I10").Value = 0.03147874789838274
$ws.Range("J10").Value = 0.03147874789838274
$ws.Range("M10").Value = 18.66868666666667
$ws.Range("N10").Value = 56.00606
$ws.Range("O10").Value = 0.2102711359030499
$ws.Range("P10").Value = 0.2102711359030499
$ws.Range("Q10").Value = 62.79631561204222
$ws.Range("R10").Value = 565.16684050838
$ws.Range("S10").Value = 0.006619072077398685
$ws.Range("T10").Value = 0.006619072077398684

$ws.Range("G11").Value = 3.363724333333333
$ws.Range("H11").Value = 10.091173
$ws.Range("I11").Value = 0.03147874789838274
$ws.Range("J11").Value = 0.03147874789838274
$ws.Range("M11").Value = 5.641943333333334
$ws.Range("N11").Value = 16.92583
$ws.Range("O11").Value = 0.06354693581733691
$ws.Range("P11").Value = 0.06354693581733692
$ws.Range("Q11").Value = 18.97794207762111
$ws.Range("R11").Value = 170.80147869859
$ws.Range("S11").Value = 0.002000377972308657
$ws.Range("T11").Value = 0.002000377972308657

$ws.Range("G12").Value = 3.363724333333333
$ws.Range("H12").Value = 10.091173
$ws.Range("I12").Value = 0.03147874789838274
$ws.Range("J12").Value = 0.03147874789838274
$ws.Range("M12").Value = 15.42507366666667
$ws.Range("N12").Value = 46.275221
$ws.Range("O12").Value = 0.1737373292074942
$ws.Range("P12").Value = 0.1737373292074942
$ws.Range("Q12").Value = 51.88569563602589
$ws.Range("R12").Value = 466.971260724233
$ws.Range("S12").Value = 0.005469033586661037
$ws.Range("T12").Value = 0.005469033586661037

$ws.Range("G13").Value = 3.363724333333333
$ws.Range("H13").Value = 10.091173
$ws.Range("I13").Value = 0.03147874789838274
$ws.Range("J13").Value = 0.03147874789838274
$ws.Range("M13").Value = 13.804248
$ws.Range("N13").Value = 41.412744
$ws.Range("O13").Value = 0.1554814732859661
$ws.Range("P13").Value = 0.1554814732859661
$ws.Range("Q13").Value = 46.433684900968
$ws.Range("R13").Value = 417.903164108712
$ws.Range("S13").Value = 0.004894362100438059
$ws.Range("T13").Value = 0.004894362100438057

$ws.Range("G14").Value = 101.145495
$ws.Range("H14").Value = 303.436485
$ws.Range("I14").Value = 0.9465500804006033
$ws.Range("J14").Value = 0.9465500804006032
$ws.Range("M14").Value = 15.090721
$ws.Range("N14").Value = 45.272163
$ws.Range("O14").Value = 0.169971412714946
$ws.Range("P14").Value = 0.169971412714946
$ws.Range("Q14").Value = 1526.358445451895
$ws.Range("R14").Value = 13737.22600906706
$ws.Range("S14").Value = 0.1608864543711363
$ws.Range("T14").Value = 0.1608864543711362

$ws.Range("G15").Value = 101.145495
$ws.Range("H15").Value = 303.436485
$ws.Range("I15").Value = 0.9465500804006033
$ws.Range("J15").Value = 0.9465500804006032
$ws.Range("N15").Value = 60.45961299999999
$ws.Range("O15").Value = 0.226991713071207
$ws.Range("P15").Value = 0.226991713071207
$ws.Range("Q15").Value = 2038.405828131145
$ws.Range("R15").Value = 18345.6524531803
$ws.Range("S15").Value = 0.2148590242578217
$ws.Range("T15").Value = 0.2148590242578216

$ws.Range("G16").Value = 101.145495
$ws.Range("H16").Value = 303.436485
$ws.Range("I16").Value = 0.9465500804006033
$ws.Range("J16").Value = 0.9465500804006032
$ws.Range("M16").Value = 18.66868666666667
$ws.Range("N16").Value = 56.00606
$ws.Range("O16").Value = 0.2102711359030499
$ws.Range("P16").Value = 0.2102711359030499
$ws.Range("Q16").Value = 1888.2535538999
$ws.Range("R16").Value = 16994.2819850991
$ws.Range("S16").Value = 0.1990321605949581
$ws.Range("T16").Value = 0.1990321605949581

$ws.Range("G17").Value = 101.145495
$ws.Range("H17").Value = 303.436485
$ws.Range("I17").Value = 0.9465500804006033
$ws.Range("J17").Value = 0.9465500804006032
$ws.Range("M17").Value = 5.641943333333334
$ws.Range("N17").Value = 16.92583
$ws.Range("O17").Value = 0.06354693581733691
$ws.Range("P17").Value = 0.06354693581733692
$ws.Range("Q17").Value = 570.65715121195
$ws.Range("R17").Value = 5135.914360907551
$ws.Range("S17").Value = 0.06015035720711222
$ws.Range("T17").Value = 0.06015035720711223

$ws.Range("G18").Value = 101.145495
$ws.Range("H18").Value = 303.436485
$ws.Range("I18").Value = 0.9465500804006033
$ws.Range("J18").Value = 0.9465500804006032
$ws.Range("M18").Value = 15.42507366666667
$ws.Range("N18").Value = 46.275221
$ws.Range("O18").Value = 0.1737373292074942
$ws.Range("P18").Value = 0.1737373292074942
$ws.Range("Q18").Value = 1560.176711426465
$ws.Range("R18").Value = 14041.59040283819
$ws.Range("S18").Value = 0.1644510829299397
$ws.Range("T18").Value = 0.1644510829299397

$ws.Range("G19").Value = 101.145495
$ws.Range("H19").Value = 303.436485
$ws.Range("I19").Value = 0.9465500804006033
$ws.Range("J19").Value = 0.9465500804006032
$ws.Range("M19").Value = 13.804248
$ws.Range("N19").Value = 41.412744
$ws.Range("O19").Value = 0.1554814732859661
$ws.Range("P19").Value = 0.1554814732859661
$ws.Range("Q19").Value = 1396.23749706276
$ws.Range("R19").Value = 12566.13747356484
$ws.Range("S19").Value = 0.1471710010396355
$ws.Range("T19").Value = 0.1471710010396355

$ws.Range("G20").Value = 2.152479666666667
$ws.Range("H20").Value = 6.457439
$ws.Range("I20").Value = 0.02014355460462176
$ws.Range("J20").Value = 0.02014355460462176
$ws.Range("M20").Value = 15.090721
$ws.Range("N20").Value = 45.272163
$ws.Range("O20").Value = 0.169971412714946
$ws.Range("P20").Value = 0.169971412714946
$ws.Range("Q20").Value = 32.48247010783967
$ws.Range("R20").Value = 292.342230970557
$ws.Range("S20").Value = 0.003423828433248216
$ws.Range("T20").Value = 0.003423828433248216

$ws.Range("G21").Value = 2.152479666666667
$ws.Range("H21").Value = 6.457439
$ws.Range("I21").Value = 0.02014355460462176
$ws.Range("J21").Value = 0.02014355460462176
$ws.Range("N21").Value = 60.45961299999999
$ws.Range("O21").Value = 0.226991713071207
$ws.Range("P21").Value = 0.226991713071207
$ws.Range("Q21").Value = 43.37936254567855
$ws.Range("R21").Value = 390.414262911107
$ws.Range("S21").Value = 0.004572419967046493
$ws.Range("T21").Value = 0.004572419967046492

$ws.Range("G22").Value = 2.152479666666667
$ws.Range("H22").Value = 6.457439
$ws.Range("I22").Value = 0.02014355460462176
$ws.Range("J22").Value = 0.02014355460462176
$ws.Range("M22").Value = 18.66868666666667
$ws.Range("N22").Value = 56.00606
$ws.Range("O22").Value = 0.2102711359030499
$ws.Range("P22").Value = 0.2102711359030499
$ws.Range("Q22").Value = 40.18396845337112
$ws.Range("R22").Value = 361.65571608034
$ws.Range("S22").Value = 0.004235608107838929
$ws.Range("T22").Value = 0.004235608107838929

$ws.Range("G23").Value = 2.152479666666667
$ws.Range("H23").Value = 6.457439
$ws.Range("I23").Value = 0.02014355460462176
$ws.Range("J23").Value = 0.02014355460462176
$ws.Range("M23").Value = 5.641943333333334
$ws.Range("N23").Value = 16.92583
$ws.Range("O23").Value = 0.06354693581733691
$ws.Range("P23").Value = 0.06354693581733692
$ws.Range("Q23").Value = 12.14416830548556
$ws.Range("R23").Value = 109.29751474937
$ws.Range("S23").Value = 0.00128006117159292
$ws.Range("T23").Value = 0.00128006117159292

$ws.Range("G24").Value = 2.152479666666667
$ws.Range("H24").Value = 6.457439
$ws.Range("I24").Value = 0.02014355460462176
$ws.Range("J24").Value = 0.02014355460462176
$ws.Range("M24").Value = 15.42507366666667
$ws.Range("N24").Value = 46.275221
$ws.Range("O24").Value = 0.1737373292074942
$ws.Range("P24").Value = 0.1737373292074942
$ws.Range("Q24").Value = 33.20215742433545
$ws.Range("R24").Value = 298.819416819019
$ws.Range("S24").Value = 0.003499687377752306
$ws.Range("T24").Value = 0.003499687377752306

$ws.Range("G25").Value = 2.152479666666667
$ws.Range("H25").Value = 6.457439
$ws.Range("I25").Value = 0.02014355460462176
$ws.Range("J25").Value = 0.02014355460462176
$ws.Range("M25").Value = 13.804248
$ws.Range("N25").Value = 41.412744
$ws.Range("O25").Value = 0.1554814732859661
$ws.Range("P25").Value = 0.1554814732859661
$ws.Range("Q25").Value = 29.713363133624
$ws.Range("R25").Value = 267.420268202616
$ws.Range("S25").Value = 0.003131949547142898
$ws.Range("T25").Value = 0.003131949547142898
